$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clases")
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Día"
$ws.Range("C1").Value = "Horario de inicio"
$ws.Range("D1").Value = "Horario de fin"
$ws.Range("E1").Value = "Cantidad de alumnos"
$ws.Range("F1").Value = "Equipamiento necesario"
$ws.Range("G1").Value = "Edificio preferencial"
$ws.Range("H1").Value = "aula_asignada"
$ws.Range("I1").Value = "Carrera"
$ws.Range("A2").Value = "Mate I comision 1"
Write-Host "done"
